$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rankings")

$ws.Range("B1").Value = "Bottom Ranking"
$ws.Range("C1").Value = "Top Ranking"
